$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $cell = $ws.Range($cellAddr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = '67.215.55'
$ws.Range("E2").Value = '  +7.86%  '
$ws.Range("D3").Value = '3.524.44'
$ws.Range("E3").Value = '  +11.37%  '
Set-TextValue "D4" '0.998'
$ws.Range("E4").Value = '  -0.17%  '
Set-TextValue "D5" '191.29'
$ws.Range("E5").Value = '  +12.58%  '
Set-TextValue "D6" '552.97'
$ws.Range("E6").Value = '  +5.65%  '
$ws.Range("D7").Value = '3.508.04'
$ws.Range("E7").Value = '  +11.06%  '
Set-TextValue "D8" '0.608'
$ws.Range("E8").Value = '  +3.01%  '
Set-TextValue "D9" '0.999'
$ws.Range("E9").Value = '  -0.13%  '
Set-TextValue "D10" '0.634'
$ws.Range("E10").Value = '  +5.89%  '
Set-TextValue "D11" '0.152'
$ws.Range("E11").Value = '  +17.06%  '
Set-TextValue "D12" '55.15'
$ws.Range("E12").Value = '  +5.27%  '
Set-TextValue "D13" '0.0000270'
$ws.Range("E13").Value = '  +8.72%  '
Set-TextValue "D14" '9.36'
$ws.Range("E14").Value = '  +4.15%  '
$ws.Range("D15").Value = '4.070.43'
$ws.Range("E15").Value = '  +11.12%  '
$ws.Range("D16").Value = '3.509.78'
$ws.Range("E16").Value = '  +11.36%  '
$ws.Range("E17").Value = '  +4.26%  '
$ws.Range("D18").Value = '67.059.30'
$ws.Range("E18").Value = '  +8.03%  '
Set-TextValue "D19" '18.21'
$ws.Range("E19").Value = '  +6.71%  '
Set-TextValue "D20" '11.93'
$ws.Range("E20").Value = '  +9.67%  '
Set-TextValue "D21" '0.996'
$ws.Range("E21").Value = '  +3.49%  '
Set-TextValue "D22" '427.93'
$ws.Range("E22").Value = '  +18.58%  '
Set-TextValue "D23" '3.91'
$ws.Range("E23").Value = '  +5.56%  '
Set-TextValue "D24" '84.84'
$ws.Range("E24").Value = '  +5.67%  '
Set-TextValue "D25" '4.18'
$ws.Range("E25").Value = '  +7.56%  '
Set-TextValue "D26" '11.20'
$ws.Range("E26").Value = '  +0.22%  '
Set-TextValue "D27" '2.91'
$ws.Range("E27").Value = '  +11.52%  '
Set-TextValue "D28" '11.99'
$ws.Range("E28").Value = '  +6.95%  '
Set-TextValue "D29" '8.98'
$ws.Range("E29").Value = '  +10.86%  '
Set-TextValue "D30" '30.32'
$ws.Range("E30").Value = '  +8.23%  '
Set-TextValue "D31" '650.00'
$ws.Range("E31").Value = '  +2.13%  '
Set-TextValue "D32" '6.69'
$ws.Range("E32").Value = '  +4.76%  '
Set-TextValue "D33" '11.72'
$ws.Range("E33").Value = '  +4.29%  '
Set-TextValue "D34" '0.111'
$ws.Range("E34").Value = '  +6.49%  '
Set-TextValue "D35" '59.29'
$ws.Range("E35").Value = '  +5.83%  '
Set-TextValue "D36" '38.71'
$ws.Range("E36").Value = '  +5.26%  '
$ws.Range("D37").Value = '0.0₃0819'
$ws.Range("E37").Value = '  +17.59%  '
Set-TextValue "D38" '1.00'
$ws.Range("E38").Value = '  -0.03%  '
Set-TextValue "D39" '0.390'
Set-TextValue "D40" '0.142'
$ws.Range("E40").Value = '  +14.69%  '
Set-TextValue "D41" '3.32'
$ws.Range("E41").Value = '  +14.26%  '
Set-TextValue "D42" '1.00'
$ws.Range("E42").Value = '  +0.30%  '
$ws.Range("D43").Value = '3.021.41'
$ws.Range("E43").Value = '  +5.04%  '
$ws.Range("E44").Value = '  +5.56%  '
Set-TextValue "D45" '2.88'
$ws.Range("E45").Value = '  +13.14%  '
Set-TextValue "D46" '3.37'
$ws.Range("E46").Value = '  +14.03%  '
Set-TextValue "D47" '0.0418'
$ws.Range("E47").Value = '  +7.50%  '
Set-TextValue "D48" '2.77'
$ws.Range("E48").Value = '  +5.37%  '
Set-TextValue "D49" '0.131'
$ws.Range("E49").Value = '  +6.65%  '
$ws.Range("E50").Value = '  +15.16%  '
Set-TextValue "D51" '140.52'
$ws.Range("E51").Value = '  +4.99%  '
